$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# README sheet (sheet1 / "README") restructuring:
# Add a title-page style block (paper title, "for the paper", "by", authors,
# "submitted to", journal name) above the existing "Data Template" section,
# without disturbing the existing Contents/data_wf/wf_/data_lamp/lamp_ blocks
# (they simply shift down by 7 rows).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("README")

# Unmerge existing merged description cells before shuffling rows around.
$ws.Range("A7:A8").UnMerge()
$ws.Range("A10:A11").UnMerge()
$ws.Range("A13:A14").UnMerge()
$ws.Range("A16:A17").UnMerge()

# Insert 7 fresh blank rows above the old row 3 ("Data Template"); this
# pushes the old rows 3-25 down to rows 10-32 and leaves rows 1 and 2
# (old title / "Supplementary Materials") untouched for now.
$ws.Range("A3:A9").EntireRow.Insert()

# Move the paper title down into row 4, and blank out the old row 1 (the
# sheet's used range will now start at row 2).
$ws.Range("A4").Value = $ws.Range("A1").Text
$ws.Range("A1").Clear()

# New title-block rows.
$ws.Range("A3").Value = "for the paper"
$ws.Range("A5").Value = "by"
$ws.Range("A6").Value = "Matthew Kuperus Heun, Zeke Marshall, Emmanuel Aramendia, and Paul E. Brockway"
$ws.Range("A7").Value = "submitted to"
$ws.Range("A8").Value = "Energies"
$ws.Range("A9").Value = ""

# Give the whole title block (rows 3-9) the same look as the existing
# header rows (row 2's shaded/bold style).
$ws.Range("A2").Copy()
$ws.Range("A3:A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-create the merged description cells at their new (shifted) locations.
$ws.Range("A14:A15").Merge()
$ws.Range("A17:A18").Merge()
$ws.Range("A20:A21").Merge()
$ws.Range("A23:A24").Merge()

# Match the (very slightly) taller wrapped-text row for the template blurb.
$ws.Rows.Item(11).RowHeight = 43.25

# Update the active selection to match the edited workbook.
$ws.Range("A7").Select()

Write-Output "README restructuring complete"
